$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the SamplesTab query (row 3, column B): the "Tumor" column should be
# derived from samp.sample_tumor_status instead of the collected `tumor` alias.
$ws.Range("B3").Value = "MATCH (s:study)<--(p:participant)<--(samp:sample)`nWHERE s.study_name in [`"LCCC 1108: Development of a Tumor Molecular Analyses Program and Its Use to Support Treatment Decisions (UNCseqTM)`"]`nWITH p,s,samp,COLLECT(DISTINCT samp.sample_tumor_status) as tumor`nRETURN  `n coalesce(samp.sample_id, '') as ``Sample ID``,`n coalesce(p.participant_id,'') as ``Participant ID``,`n coalesce(s.study_name, '') as ``Study Name``,`n coalesce(s.phs_accession,'') as ``Accession``,`ncoalesce(samp.sample_tumor_status,'') as ``Tumor``,`ncoalesce(samp.sample_type,'') as ``Analyte Type```nORDER By samp.sample_id LIMIT 100"

# Update the selected cell in the sheet view from E13 to C13.
$ws.Range("C13").Select()
